$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml): update F2 130 -> 132 and F3 20 -> 21
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 132
$ws1.Range("F3").Value = 21

# Sheet "全部类型" (sheet4.xml): update F2 130 -> 132 and F3 20 -> 21
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 132
$ws4.Range("F3").Value = 21
